# "Fix titles to be easily understandable"
# The three calculated-column headers on the "Entropie" sheet used a plain
# hyphen ("Title - explanation") to separate the short title from the
# explanation of how the value is computed. That's replaced with a colon
# ("Title: explanation") so the relationship reads more clearly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entropie")

$ws.Range("E13").Value = "Entropie pro Zeichen: log_2(Mögliche Zeichen)"
$ws.Range("G13").Value = "Entropie Gesamt: Anzahl Zeichen*Entropie"
$ws.Range("I13").Value = "Anzahl Kombinationen: Mögliche Zeichen^Anzahl"

# Column B was nudged a bit wider (title column now holds slightly longer
# labels such as "Zeichensatz B ").
$ws.Columns.Item(2).ColumnWidth = 13.25

# Leave the cursor resting on the "Anzahl Kombinationen" result cell, as in
# the saved workbook.
$ws.Range("I14").Select()
